# Update the Wnt7a-Fzd10 LR-pairs sheet with newly computed TPM-based values.
# - Row 2 (FAPs -> ECs) becomes (FAPs -> FAPs) with recalculated metrics
# - Row 3 (FAPs -> FAPs) becomes (FAPs -> MuSCs) with recalculated metrics
# - Rows 4 and 5 (FAPs -> MuSCs, FAPs -> Resolving-Mac) are removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two trailing rows (old rows 4 and 5) first so the remaining
# two data rows (2 and 3) can be updated in place.
$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(4).Delete() | Out-Null

# --- Row 2: target cluster changes from "ECs" to "FAPs", metrics updated ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.085107
$ws.Range("N2").Value = 0.255321
$ws.Range("O2").Value = 0.803017436594203
$ws.Range("P2").Value = 0.8030174365942029
$ws.Range("Q2").Value = 0.017917179544
$ws.Range("R2").Value = 0.161254615896
$ws.Range("S2").Value = 0.803017436594203
$ws.Range("T2").Value = 0.8030174365942029

# --- Row 3: target cluster changes from "FAPs" to "MuSCs", metrics updated ---
$ws.Range("D3").Value = "MuSCs"
$ws.Range("M3").Value = 0.020877
$ws.Range("N3").Value = 0.06263099999999999
$ws.Range("O3").Value = 0.1969825634057971
$ws.Range("P3").Value = 0.1969825634057971
$ws.Range("Q3").Value = 0.004395137384
$ws.Range("R3").Value = 0.039556236456
$ws.Range("S3").Value = 0.1969825634057971
$ws.Range("T3").Value = 0.1969825634057971

$wb.Save()
